$d = $word.ActiveDocument
$d.Content.Find.Execute(".: `${", $true, $false, $false, $false, $false, $true, 1, $false, ".: Φ.15.1/`${", 2)
